# Insert a new weekly price record at row 92 of "Sheet1", pushing the
# existing rows 92-136 down to 93-137 (preserving their data/format).
# The new row holds the same Mercado/Región/Codreg/Categoría/Calidad/
# Clasificación values as the surrounding "Ají" records, with fresh
# Fecha / Variedad / Volumen / Precio / Unidad / Origen figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 92:136 down by one, duplicating row 92's formatting
# (this is what gives the new D92 its date number format).
$ws.Rows("92:92").Insert()

$ws.Range("A92").Value = 11
$ws.Range("B92").Value = "Vega Monumental Concepción"
$ws.Range("C92").Value = "Bíobío"
$ws.Range("D92").Value = 44839
$ws.Range("E92").Value = 8
$ws.Range("F92").Value = 100112021
$ws.Range("G92").Value = "Ají"
$ws.Range("H92").Value = "Americana (o)"
$ws.Range("I92").Value = "Primera"
$ws.Range("J92").Value = 100
$ws.Range("K92").Value = 80000
$ws.Range("L92").Value = 82000
$ws.Range("M92").Value = 81000
$ws.Range("N92").Value = "`$/caja 25 kilos"
$ws.Range("O92").Value = "Provincia de Limarí"
$ws.Range("P92").Value = 3240
$ws.Range("Q92").Value = 25
$ws.Range("R92").Value = "Hortaliza"
